$d = $word.ActiveDocument

# Replace "Killing the animal to " with "Killing the animal "
$d.Content.Find.Execute("Killing the animal to ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Killing the animal ", 2)

# Replace "be molded" with "for molding"
$d.Content.Find.Execute("be molded", $true, $false, $false, $false, $false,
                         $true, 1, $false, "for molding", 2)
